# Add the new "Stability Test Report" row to the status lookup table and
# widen column C to fit the new, longer text (mirrors the manual edit:
# "add static and fix accept request").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "STAB"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "Stability Test Report"
$ws.Range("D4").Value = $true

# Column C needs to grow from 14.86 to fit "Stability Test Report" (~21 chars).
$ws.Range("C1").ColumnWidth = 20.1666666666666668

# Leave the selection where the author's last edit left it.
$ws.Range("C10").Select() | Out-Null
